$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 133 ("Especial" quality,
# fecha 2023-08-28 / serial 45166), pushing the former rows 133-148 down
# to rows 134-149 unchanged.
$ws.Rows.Item(133).Insert()

$ws.Cells.Item(133, 1).Value = 5
$ws.Cells.Item(133, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(133, 3).Value = "Maule"
$ws.Cells.Item(133, 4).Value = 45166
$ws.Cells.Item(133, 5).Value = 7
$ws.Cells.Item(133, 6).Value = "Fruta"
$ws.Cells.Item(133, 7).Value = 100107
$ws.Cells.Item(133, 8).Value = "Otros"
$ws.Cells.Item(133, 9).Value = 100107002
$ws.Cells.Item(133, 10).Value = "Chirimoya"
$ws.Cells.Item(133, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(133, 12).Value = "Especial"
$ws.Cells.Item(133, 13).Value = 100
$ws.Cells.Item(133, 14).Value = 30000
$ws.Cells.Item(133, 15).Value = 30000
$ws.Cells.Item(133, 16).Value = 30000
$ws.Cells.Item(133, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(133, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(133, 19).Value = 3000
$ws.Cells.Item(133, 20).Value = 10
